# Kanban Board update: move "Account System: Activity Log" from the
# "Not Started" column (A9) to the "Done" column (append at C16).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kanban Boad")

# Remove the task from the "Not Started" column (A9)
$ws.Range("A9").ClearContents()

# Append the task to the bottom of the "Done" column (C16)
$ws.Range("C16").Value = "Account System: Activity Log"

# Force recalculation so the COUNTA/percentage formulas update
$excel.Calculate()

# Update the selected cell to reflect where the user was last working
$ws.Range("B13").Select()

$wb.Save()
